$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Cell D122: "1A" -> "5"
$ws.Range("D122").NumberFormat = "@"
$ws.Range("D122").Value = "5"
$ws.Range("D122").Style = $ws.Range("A122").Style

# 2) Insert a new row at 221 (fills the gap for código 221), shifting rows 221-250 down to 222-251
$ws.Rows(221).Insert()

$ws.Range("A221").Value = 221
$ws.Range("B221").NumberFormat = "@"
$ws.Range("B221").Value = "1A"
$ws.Range("C221").NumberFormat = "@"
$ws.Range("C221").Value = "1A"
$ws.Range("D221").NumberFormat = "@"
$ws.Range("D221").Value = "5"
$ws.Range("E221").NumberFormat = "@"
$ws.Range("E221").Value = "1A"

# Match the formatting/style of the surrounding data rows (no explicit style)
$ws.Range("A221:E221").Style = $ws.Range("A220:E220").Style

# 3) Cell D246 (previously row 245 before the insertion shifted it down): "5" -> "1A"
$ws.Range("D246").NumberFormat = "@"
$ws.Range("D246").Value = "1A"
$ws.Range("D246").Style = $ws.Range("A246").Style
